$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in missing "Original" (column C) values for the existing EPFL rows ---
$ws.Range("C6").Value = 26.2822
$ws.Range("C7").Value = 26.1257
$ws.Range("C8").Value = 26.4489

# --- Additional EPFL dataset rows (9-11) ---
$ws.Range("A9").Value = "Magnets_1"
$ws.Range("B9").Value = "EPFL"

$ws.Range("A10").Value = "Stone_Pillars_Outside"
$ws.Range("B10").Value = "EPFL"

$ws.Range("A11").Value = "Vespa"
$ws.Range("B11").Value = "EPFL"

$ws.Range("C9").Value = 25.7569
$ws.Range("C10").Value = 26.1003
$ws.Range("C11").Value = 26.5211

# --- STANFORD dataset rows (12-20) ---
$ws.Range("A12").Value = "Amethyst"
$ws.Range("A13").Value = "Bracelet"
$ws.Range("A14").Value = "Chess"
$ws.Range("A15").Value = "Eucalyptus Flowers"

$ws.Range("B12").Value = "STANFORD"

$ws.Range("A16").Value = "Jelly Beens"
$ws.Range("A17").Value = "Lego Bulldozer"
$ws.Range("A18").Value = "Lego Knights"
$ws.Range("A19").Value = "Lego Truck"
$ws.Range("A20").Value = "Treasure"

$ws.Range("B13").Value = "STANFORD"
$ws.Range("B14").Value = "STANFORD"
$ws.Range("B15").Value = "STANFORD"
$ws.Range("B16").Value = "STANFORD"
$ws.Range("B17").Value = "STANFORD"
$ws.Range("B18").Value = "STANFORD"
$ws.Range("B19").Value = "STANFORD"
$ws.Range("B20").Value = "STANFORD"

$ws.Range("C12").Value = 26.7365
$ws.Range("C13").Value = 25.6613
$ws.Range("C14").Value = 25.9557
$ws.Range("C15").Value = 26.9736
$ws.Range("C16").Value = 25.6778
$ws.Range("C17").Value = 26.433
$ws.Range("C18").Value = 25.7786
$ws.Range("C19").Value = 26.1476
$ws.Range("C20").Value = 26.6825

# --- Summary row 21: Average EPFL ---
$ws.Range("A21").Value = "Average EPFL"
$ws.Range("A21").Font.Bold = $true
$ws.Range("A21").IndentLevel = 0
$ws.Range("A21:B21").MergeCells = $true
$ws.Range("B21").Font.Bold = $false
$ws.Range("B21").IndentLevel = 0
$ws.Range("C21").Formula = "=AVERAGE(C2:C11)"
$ws.Range("D21:F21").Formula = "=AVERAGE(D2:D11)"

# --- Summary row 22: Average STANFORD ---
$ws.Range("A22").Value = "Average STANFORD"
$ws.Range("A22").Font.Bold = $true
$ws.Range("A22").IndentLevel = 0
$ws.Range("A22:B22").MergeCells = $true
$ws.Range("C22").Formula = "=AVERAGE(C12:C20)"
$ws.Range("D22:F22").Formula = "=AVERAGE(D12:D20)"

# --- Summary row 23: Average Overall ---
$ws.Range("A23").Value = "Average Overall"
$ws.Range("A23").Font.Bold = $true
$ws.Range("A23").IndentLevel = 0
$ws.Range("A23:B23").MergeCells = $true
$ws.Range("C23").Formula = "=AVERAGE(C2:C20)"
$ws.Range("D23:F23").Formula = "=AVERAGE(D2:D20)"

# --- Restore the cursor position as recorded in the saved workbook ---
$ws.Range("G29").Select()
